$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Local" input value in the Best/Test table (C3)
$ws.Range("C3").Value = 3402436

# Update the "Local Test" table inputs (F2 = Starting, F3 = Test)
$ws.Range("F2").Value = 3402436
$ws.Range("F3").Value = 3115976

# Update the "Scoreboard" table input (C10 = Local)
$ws.Range("C10").Value = 12121765

# Update the active selection to F3 as in the diff
$ws.Range("F3").Select()
